# wg() only takes 'df' and 'tb'
# Update the AJ ("Wohngeld" max transfer) and AQ (min/max income) formulas
# for every household row (2..13) to match the simplified wg() signature,
# and move the sheet view / selection to show the new formulas (column Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $ws.Range("AJ$r").Formula = "=MAX((1-AF$r)*(AH$r-AI$r),0)"
    $ws.Range("AQ$r").Formula = "=MIN(MAX(I$r,AO$r),AN$r)"
}

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 17
$ws.Range("Q4").Select()
